$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so values like "591.36" or
# "0.0000237" are not auto-converted into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '62.308.05'
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("D3").Value = '3.109.95'
$ws.Range("E3").Value = '  -2.00%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '591.36'
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("D6").Value = '132.78'
$ws.Range("E6").Value = '  -1.18%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '3.108.50'
$ws.Range("E8").Value = '  -1.92%  '
$ws.Range("D9").Value = '0.504'
$ws.Range("E9").Value = '  +0.34%  '
$ws.Range("E10").Value = '  -0.37%  '
$ws.Range("D11").Value = '5.33'
$ws.Range("E11").Value = '  +2.25%  '
$ws.Range("D12").Value = '0.445'
$ws.Range("E12").Value = '  -1.00%  '
$ws.Range("D13").Value = '0.0000237'
$ws.Range("E13").Value = '  +1.75%  '
$ws.Range("D14").Value = '33.98'
$ws.Range("E14").Value = '  +3.02%  '
$ws.Range("D15").Value = '3.624.47'
$ws.Range("E15").Value = '  -1.99%  '
$ws.Range("E16").Value = '  +1.09%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '62.677.93'
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.120.07'
$ws.Range("E18").Value = '  -1.63%  '
$ws.Range("D19").Value = '6.44'
$ws.Range("E19").Value = '  -1.38%  '
$ws.Range("D20").Value = '453.41'
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("D21").Value = '13.69'
$ws.Range("E21").Value = '  -0.85%  '
$ws.Range("D22").Value = '0.683'
$ws.Range("E22").Value = '  -2.57%  '
$ws.Range("D23").Value = '7.49'
$ws.Range("E23").Value = '  -1.61%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '81.70'
$ws.Range("E24").Value = '  -0.57%  '
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = '12.98'
$ws.Range("E25").Value = '  -2.71%  '
$ws.Range("D26").Value = '0.997'
$ws.Range("E26").Value = '  -0.30%  '
$ws.Range("E27").Value = '  +0.34%  '
$ws.Range("D28").Value = '2.63'
$ws.Range("E28").Value = '  -1.53%  '
$ws.Range("D29").Value = '2.04'
$ws.Range("E29").Value = '  +1.74%  '
$ws.Range("D30").Value = '7.55'
$ws.Range("E30").Value = '  -3.13%  '
$ws.Range("D31").Value = '6.58'
$ws.Range("E31").Value = '  -4.68%  '
$ws.Range("D32").Value = '26.51'
$ws.Range("E32").Value = '  -2.39%  '
$ws.Range("D33").Value = '0.0990'
$ws.Range("E33").Value = '  -2.84%  '
$ws.Range("D34").Value = '2.35'
$ws.Range("E34").Value = '  -1.40%  '
$ws.Range("D35").Value = '0.997'
$ws.Range("E35").Value = '  -3.19%  '
$ws.Range("D36").Value = '5.78'
$ws.Range("E36").Value = '  +0.16%  '
$ws.Range("D37").Value = '50.74'
$ws.Range("E37").Value = '  -0.67%  '
$ws.Range("D38").Value = '0.0₃0707'
$ws.Range("E38").Value = '  +2.61%  '
$ws.Range("D39").Value = '0.0382'
$ws.Range("E39").Value = '  -0.40%  '
$ws.Range("D40").Value = '8.01'
$ws.Range("E40").Value = '  +0.57%  '
$ws.Range("D41").Value = '0.110'
$ws.Range("E41").Value = '  -1.51%  '
$ws.Range("E42").Value = '  -2.12%  '
$ws.Range("D43").Value = '382.72'
$ws.Range("E43").Value = '  -6.88%  '
$ws.Range("D44").Value = '2.740.56'
$ws.Range("E44").Value = '  -6.54%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").Value = '0.245'
$ws.Range("E46").Value = '  -1.58%  '
$ws.Range("D47").Value = '125.73'
$ws.Range("E47").Value = '  +0.90%  '
$ws.Range("D48").Value = '35.03'
$ws.Range("E48").Value = '  -1.49%  '
$ws.Range("D49").Value = '2.06'
$ws.Range("E49").Value = '  -2.91%  '
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").Value = '24.46'
$ws.Range("E51").Value = '  -3.24%  '
